$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 9
$ws.Range("B9").Value = 2022
$ws.Range("C9").Value = 2023
$ws.Range("D9").Value = 2024
$ws.Range("E9").Value = 2025
$ws.Range("F9").Value = 2026
$ws.Range("G9").Value = 2027
$ws.Range("H9").Value = 2028
$ws.Range("I9").Value = 2029
$ws.Range("J9").Value = 2030
$ws.Range("K9").Value = 2031
$ws.Range("L9").Value = 2032

# Row 12
$ws.Range("B12").Value = 13.83
$ws.Range("C12").Value = 14.25
$ws.Range("D12").Value = 14.81
$ws.Range("E12").Value = 15.08
$ws.Range("F12").Value = 15.36
$ws.Range("G12").Value = 15.86
$ws.Range("H12").Value = 15.93
$ws.Range("I12").Value = 15.9
$ws.Range("J12").Value = 16
$ws.Range("K12").Value = 15.95
$ws.Range("L12").Value = 16.08

# Row 14
$ws.Range("B14").Value = 17.32
$ws.Range("C14").Value = 17.97
$ws.Range("D14").Value = 18.86
$ws.Range("E14").Value = 19.53
$ws.Range("F14").Value = 22.44
$ws.Range("G14").Value = 23.15
$ws.Range("H14").Value = 23.24
$ws.Range("I14").Value = 23.21
$ws.Range("J14").Value = 23.32
$ws.Range("K14").Value = 23.29
$ws.Range("L14").Value = 23.42

# Row 16
$ws.Range("B16").Value = 6.59
$ws.Range("C16").Value = 9.05
$ws.Range("D16").Value = 11.61
$ws.Range("E16").Value = 13.88
$ws.Range("F16").Value = 17.21
$ws.Range("G16").Value = 19.46
$ws.Range("H16").Value = 19.57
$ws.Range("I16").Value = 19.54
$ws.Range("J16").Value = 19.65
$ws.Range("K16").Value = 19.62
$ws.Range("L16").Value = 19.77

# Row 17
$ws.Range("B17").Value = 19.8
$ws.Range("C17").Value = 20.09
$ws.Range("D17").Value = 20.66
$ws.Range("E17").Value = 20.96
$ws.Range("F17").Value = 23.19
$ws.Range("G17").Value = 23.56
$ws.Range("H17").Value = 23.66
$ws.Range("I17").Value = 23.64
$ws.Range("J17").Value = 23.74
$ws.Range("K17").Value = 23.71
$ws.Range("L17").Value = 23.85

# Row 18
$ws.Range("B18").Value = 25.02
$ws.Range("C18").Value = 25.02
$ws.Range("D18").Value = 25.22
$ws.Range("E18").Value = 25.2
$ws.Range("F18").Value = 29.47
$ws.Range("G18").Value = 29.48
$ws.Range("H18").Value = 29.54
$ws.Range("I18").Value = 29.51
$ws.Range("J18").Value = 29.61
$ws.Range("K18").Value = 29.59
$ws.Range("L18").Value = 29.7

# Row 19
$ws.Range("B19").Value = -6.56
$ws.Range("C19").Value = -6.76
$ws.Range("D19").Value = -6.7
$ws.Range("E19").Value = -6.84
$ws.Range("F19").Value = -5.21
$ws.Range("G19").Value = -5.26
$ws.Range("H19").Value = -5.12
$ws.Range("I19").Value = -5.15
$ws.Range("J19").Value = -5.02
$ws.Range("K19").Value = -5.06
$ws.Range("L19").Value = -4.88

# Row 20
$ws.Range("B20").Value = 7.87
$ws.Range("C20").Value = 12.05
$ws.Range("D20").Value = 16.1
$ws.Range("E20").Value = 19.63
$ws.Range("F20").Value = 23.86
$ws.Range("G20").Value = 26.91
$ws.Range("H20").Value = 27.03
$ws.Range("I20").Value = 27.01
$ws.Range("J20").Value = 27.12
$ws.Range("K20").Value = 27.09
$ws.Range("L20").Value = 27.24

# Row 21
$ws.Range("B21").Value = 29.84
$ws.Range("C21").Value = 29.75
$ws.Range("D21").Value = 29.85
$ws.Range("E21").Value = 29.82
$ws.Range("F21").Value = 32.7
$ws.Range("G21").Value = 32.68
$ws.Range("H21").Value = 32.75
$ws.Range("I21").Value = 32.71
$ws.Range("J21").Value = 32.79
$ws.Range("K21").Value = 32.77
$ws.Range("L21").Value = 32.86

# Row 23
$ws.Range("B23").Value = 19.18
$ws.Range("C23").Value = 19.9
$ws.Range("D23").Value = 20.74
$ws.Range("E23").Value = 21.46
$ws.Range("F23").Value = 24.47
$ws.Range("G23").Value = 25.17
$ws.Range("H23").Value = 25.22
$ws.Range("I23").Value = 25.18
$ws.Range("J23").Value = 25.23
$ws.Range("K23").Value = 25.21
$ws.Range("L23").Value = 25.26

# Row 24
$ws.Range("B24").Value = 7.7
$ws.Range("C24").Value = 8.05
$ws.Range("D24").Value = 9.17
$ws.Range("E24").Value = 9.47
$ws.Range("F24").Value = 11.32
$ws.Range("G24").Value = 12.06
$ws.Range("H24").Value = 12.43
$ws.Range("I24").Value = 12.45
$ws.Range("J24").Value = 12.88
$ws.Range("K24").Value = 12.81
$ws.Range("L24").Value = 13.4

# Row 26
$ws.Range("B26").Value = 15.62
$ws.Range("C26").Value = 16.32
$ws.Range("D26").Value = 17.26
$ws.Range("E26").Value = 17.98
$ws.Range("F26").Value = 19.44
$ws.Range("G26").Value = 20.19
$ws.Range("H26").Value = 20.32
$ws.Range("I26").Value = 20.29
$ws.Range("J26").Value = 20.4
$ws.Range("K26").Value = 20.36
$ws.Range("L26").Value = 20.51

# Row 27
$ws.Range("B27").Value = 21.04
$ws.Range("C27").Value = 21.6
$ws.Range("D27").Value = 22.38
$ws.Range("E27").Value = 22.91
$ws.Range("F27").Value = 28.65
$ws.Range("G27").Value = 29.27
$ws.Range("H27").Value = 29.31
$ws.Range("I27").Value = 29.27
$ws.Range("J27").Value = 29.37
$ws.Range("K27").Value = 29.36
$ws.Range("L27").Value = 29.45

# Row 29
$ws.Range("B29").Value = 6.71
$ws.Range("C29").Value = 6.58
$ws.Range("D29").Value = 6.35
$ws.Range("E29").Value = 5.66
$ws.Range("F29").Value = -1.42
$ws.Range("G29").Value = -1.61
$ws.Range("H29").Value = -1.6
$ws.Range("I29").Value = -1.65
$ws.Range("J29").Value = -1.58
$ws.Range("K29").Value = -1.68
$ws.Range("L29").Value = -1.58

# Row 31
$ws.Range("B31").Value = -0.26
$ws.Range("C31").Value = -0.23
$ws.Range("D31").Value = -0.24
$ws.Range("E31").Value = -0.26
$ws.Range("F31").Value = -3.13
$ws.Range("G31").Value = -3.17
$ws.Range("H31").Value = -3.16
$ws.Range("I31").Value = -3.16
$ws.Range("J31").Value = -3.17
$ws.Range("K31").Value = -3.16
$ws.Range("L31").Value = -3.16

# Row 32
$ws.Range("B32").Value = 25.42
$ws.Range("C32").Value = 24.82
$ws.Range("D32").Value = 24.25
$ws.Range("E32").Value = 22.07
$ws.Range("F32").Value = 4.53
$ws.Range("G32").Value = 3.82
$ws.Range("H32").Value = 3.84
$ws.Range("I32").Value = 3.63
$ws.Range("J32").Value = 3.96
$ws.Range("K32").Value = 3.49
$ws.Range("L32").Value = 3.98

# Row 35
$ws.Range("B35").Value = 14.49
$ws.Range("C35").Value = 14.77
$ws.Range("D35").Value = 15.18
$ws.Range("E35").Value = 15.28
$ws.Range("F35").Value = 15.21
$ws.Range("G35").Value = 15.55
$ws.Range("H35").Value = 15.62
$ws.Range("I35").Value = 15.58
$ws.Range("J35").Value = 15.68
$ws.Range("K35").Value = 15.63
$ws.Range("L35").Value = 15.76

# Row 36
$ws.Range("B36").Value = 19.33
$ws.Range("C36").Value = 19.83
$ws.Range("D36").Value = 20.56
$ws.Range("E36").Value = 21.07
$ws.Range("F36").Value = 24.12
$ws.Range("G36").Value = 24.67
$ws.Range("H36").Value = 24.76
$ws.Range("I36").Value = 24.73
$ws.Range("J36").Value = 24.84
$ws.Range("K36").Value = 24.81
$ws.Range("L36").Value = 24.94

# Row 37
$ws.Range("B37").Value = 27.4
$ws.Range("C37").Value = 27.34
$ws.Range("D37").Value = 27.48
$ws.Range("E37").Value = 27.43
$ws.Range("F37").Value = 30.96
$ws.Range("G37").Value = 30.94
$ws.Range("H37").Value = 31.01
$ws.Range("I37").Value = 30.98
$ws.Range("J37").Value = 31.08
$ws.Range("K37").Value = 31.06
$ws.Range("L37").Value = 31.17

# Row 38
$ws.Range("B38").Value = 6.71
$ws.Range("C38").Value = 6.58
$ws.Range("D38").Value = 6.35
$ws.Range("E38").Value = 5.66
$ws.Range("F38").Value = -1.42
$ws.Range("G38").Value = -1.61
$ws.Range("H38").Value = -1.6
$ws.Range("I38").Value = -1.65
$ws.Range("J38").Value = -1.58
$ws.Range("K38").Value = -1.68
$ws.Range("L38").Value = -1.58
